$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply same style (style index 1, quotePrefix) as existing B column cells BEFORE
# setting values, so the quote-prefix formatting carries through
$ws.Range("B5:B8").Style = $ws.Range("B3").Style

# Update existing row 4 (B4) to new value, keep A4 = PREPROD
# Leading apostrophe forces Excel to treat numeric-looking text as text (quote prefix)
$ws.Range("B4").Value = "'0420172008282  "

# Add new rows 5-8
$ws.Range("A5").Value = "PREPROD"
$ws.Range("B5").Value = "'0420172008281"

$ws.Range("A6").Value = "PREPROD"
$ws.Range("B6").Value = "'1220170301396"

$ws.Range("A7").Value = "PREPROD"
$ws.Range("B7").Value = "'1120170200917"

$ws.Range("A8").Value = "PREPROD"
$ws.Range("B8").Value = "'1220170301397"

# Update selection to B5 as seen in diff
$ws.Range("B5").Select()
